$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 (only the cells that actually change)
$ws.Range("B2").Value = 389.0

$ws.Range("A3").Value = 3.0
$ws.Range("B3").Value = 3.99

# Insert two new data rows (4 and 5) before the Total row
$ws.Range("A4").Value = 1.0
$ws.Range("B4").Value = 312.5

$ws.Range("A5").Value = 1.0
$ws.Range("B5").Value = 365.66

# Move the Total row down to row 6
$ws.Range("A6").Value = "Total"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1079.13"
